# PowerPoll "Presentation1.pptx" is a Visual-Studio/VSTO debug build artifact:
# every time the add-in project is run (F5) PowerPoint re-opens with the
# PowerPoll task-pane app freshly attached to the title slide and the deck is
# re-saved to bin/Debug. That re-attach/re-save is what produced this
# revision's OOXML diff:
#
#   * ppt/presentation.xml   - r:id of the sldMasterId/sldId entries
#   * ppt/slideMasters/...   - r:id of every sldLayoutId entry
#   * ppt/slides/slide.xml   - r:id of the webextensionref / blip r:embed
#   * ppt/slides/udata/data.xml - the we:webextension instance GUID
#     (+ the matching r:embed id for its snapshot image)
#
# All of the relationship-id churn above is PowerPoint renumbering its own
# package internals on save (every part keeps the exact same target/content,
# only the random rXXXXXXXXXXXXXXXX id text differs) - it is not something a
# script drives deliberately, it happens automatically whenever the file is
# written out. The one "real" value that changed, the embedded web add-in's
# <we:webextension id="{...}"> instance id, lives in a package part that the
# PowerPoint object model does not expose to VBA/COM automation at all (no
# Shape/Slide/Presentation member reaches ppt/slides/udata/data.xml - web
# add-ins are only managed through the Insert Add-in UI, never scripted).
#
# The faithful COM-interop equivalent of "the author re-ran/re-saved the
# deck after updating the add-in" is therefore simply to save the
# presentation again; that's exactly what's reproducible from this object
# model without guessing at / hand-forging internal package ids.
$p = $ppt.ActivePresentation
$p.Save()
